# Update cryptocurrency Price (D) and Volume(1h) (E) columns with refreshed
# market data fetched by the GitHub Actions symbol-list updater.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'299.09"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'-1.60%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'31.39"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'-1.66%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.112"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'-2.16%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.07928"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'1.06%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'2.312"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'-1.59%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'7.809"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'-2.37%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'3.863"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'-0.19%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.9229"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'1.11%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.1746"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'0.64%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.07585"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'2.78%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.09340"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'14.83%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.03006"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'-1.64%"
$ws.Range("E13").Style = "Normal"
$ws.Range("E14").Value = "'0.89%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.001504"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'-0.64%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.005841"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'-5.51%"
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'3.477"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'-0.67%"
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'2.267"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'1.24%"
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'-0.11%"
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'0.1308"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'-0.23%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'4.017"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'-13.81%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.1700"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'8.65%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.04622"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'-0.30%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.001250"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'-0.97%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.004481"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'-1.19%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.0001250"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'-7.32%"
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'0.0003395"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'23.88%"
$ws.Range("E27").Style = "Normal"
$ws.Range("D39").Value = "'0.01741"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'-2.76%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.04622"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'0.70%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.006979"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'-4.56%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.1362"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'-0.06%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.002190"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'-2.14%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.01031"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'-5.74%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.00006290"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'-2.70%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.00000000749"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'0.03%"
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.007975"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'-19.41%"
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.7465"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'-9.03%"
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.00002099"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'0.03%"
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.0001999"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'0.03%"
$ws.Range("E50").Style = "Normal"
